$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Capture the existing values (rows 86-88, the last three product rows) and
# the current running total before we shuffle anything.
# ---------------------------------------------------------------------------
$row86 = @{
    H = $ws.Range("H86").Value2
    L = $ws.Range("L86").Value2
    N = $ws.Range("N86").Value2
    P = $ws.Range("P86").Value2
    Q = $ws.Range("Q86").Value2
    C = $ws.Range("C86").Value2
}
$row87 = @{
    H = $ws.Range("H87").Value2
    L = $ws.Range("L87").Value2
    N = $ws.Range("N87").Value2
    P = $ws.Range("P87").Value2
    Q = $ws.Range("Q87").Value2
    C = $ws.Range("C87").Value2
}
$row88 = @{
    H = $ws.Range("H88").Value2
    L = $ws.Range("L88").Value2
    N = $ws.Range("N88").Value2
    P = $ws.Range("P88").Value2
    Q = $ws.Range("Q88").Value2
    C = $ws.Range("C88").Value2
}

$oldTotal = $ws.Range("P89").Value2
$oldRowHeight88 = $ws.Rows("88:88").RowHeight

# ---------------------------------------------------------------------------
# Make room: insert a fresh row just above the totals row (row 89), pushing
# the totals row down to 90 and the footer row down to 91.
# ---------------------------------------------------------------------------
$ws.Rows("89:89").Insert()

# Build the new row 89 as a copy of the data-row formatting (row 88 still has
# the original product-row styling at this point) and restore its geometry.
$ws.Range("A88:Q88").Copy()
$ws.Range("A89:Q89").PasteSpecial(-4122)
$ws.Rows("89:89").RowHeight = $oldRowHeight88

$ws.Range("A89:B89").Merge()
$ws.Range("C89:G89").Merge()
$ws.Range("H89:K89").Merge()
$ws.Range("L89:M89").Merge()
$ws.Range("N89:O89").Merge()

# ---------------------------------------------------------------------------
# Shift the three existing products down one slot (86->87, 87->88, 88->89)
# and place the brand-new product into row 86.
# ---------------------------------------------------------------------------
$ws.Range("A89").Value2 = 83
$ws.Range("C89").Value2 = $row88.C
$ws.Range("H89").Value2 = $row88.H
$ws.Range("L89").Value2 = $row88.L
$ws.Range("N89").Value2 = $row88.N
$fmt = $ws.Range("P89").NumberFormat
$ws.Range("P89").NumberFormat = "@"
$ws.Range("P89").Value2 = $row88.P
$ws.Range("P89").NumberFormat = $fmt
$ws.Range("Q89").Value2 = $row88.Q

$ws.Range("C88").Value2 = $row87.C
$ws.Range("H88").Value2 = $row87.H
$ws.Range("L88").Value2 = $row87.L
$ws.Range("N88").Value2 = $row87.N
$fmt = $ws.Range("P88").NumberFormat
$ws.Range("P88").NumberFormat = "@"
$ws.Range("P88").Value2 = $row87.P
$ws.Range("P88").NumberFormat = $fmt
$ws.Range("Q88").Value2 = $row87.Q

$ws.Range("C87").Value2 = $row86.C
$ws.Range("H87").Value2 = $row86.H
$ws.Range("L87").Value2 = $row86.L
$ws.Range("N87").Value2 = $row86.N
$fmt = $ws.Range("P87").NumberFormat
$ws.Range("P87").NumberFormat = "@"
$ws.Range("P87").Value2 = $row86.P
$ws.Range("P87").NumberFormat = $fmt
$ws.Range("Q87").Value2 = $row86.Q

$ws.Range("C86").Value2 = "فلامنجو شفرات للنساء"
$ws.Range("H86").Value2 = "4:0"
$ws.Range("L86").Value2 = "0"
$ws.Range("N86").Value2 = "40.00"
$fmt = $ws.Range("P86").NumberFormat
$ws.Range("P86").NumberFormat = "@"
$ws.Range("P86").Value2 = "40.0000"
$ws.Range("P86").NumberFormat = $fmt
$ws.Range("Q86").Value2 = "1:0"

# ---------------------------------------------------------------------------
# Update the grand-total cell (now row 90) to include the new product's
# selling price, and refresh the generated timestamp in the footer (row 91).
# ---------------------------------------------------------------------------
$ws.Range("P90").Value2 = [double]$oldTotal + 40.0

$ws.Range("A91").Value2 = "Monday, 15 September, 2025 6:55 PM"
